# Generate Report for Handoff
# The f2bb82f3-... file has been handed off: its status moves from
# "Handed back: in sync with en-US" to "Ready for handoff" on every sheet
# that references it, and the per-locale handoff datetime stamps are
# recorded for the zh-cn and de-de target files.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: row for f2bb82f3-...md (row 3) ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B3").Value = "Ready for handoff"
$overview.Range("C3").Value = "Ready for handoff"

# --- zh-cn sheet: row for f2bb82f3-...md (row 3) ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("B3").Value = "Ready for handoff"
$zhcn.Range("D3").Value = "2016-02-17 04:24:29"

# --- de-de sheet: row for f2bb82f3-...md (row 3) ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("B3").Value = "Ready for handoff"
$dede.Range("D3").Value = "2016-02-17 04:24:39"
